$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the slide-number placeholder ("Slide Number Placeholder 1") from
# slide 1 by turning off the slide-number footer for this slide. This
# drops the corresponding <p:sp> (ph type="sldNum") from the slide's
# shape tree, matching the authored edit.
$hf = $s.HeadersFooters
$hf.SlideNumber.Visible = 0
